# Fixing minor bug with timestamp on JUNC folders
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("גיליון1")

# Simple numeric value corrections (column B), labels unchanged
$ws.Range("B2").Value = 1125
$ws.Range("B3").Value = 1.347555555555556
$ws.Range("B4").Value = 1598
$ws.Range("B7").Value = 82
$ws.Range("B8").Value = 193
$ws.Range("B9").Value = 418
$ws.Range("B10").Value = 177
$ws.Range("B11").Value = 728
$ws.Range("B12").Value = 82

# Rows 14-27 get new label/value pairs (shifted JUNC folder entries),
# and a new row 28 is added below with the entry that fell off the end.
$ws.Range("A14").Value = "EEr"
$ws.Range("B14").Value = 82

$ws.Range("A15").Value = "DStl"
$ws.Range("B15").Value = 191

$ws.Range("A16").Value = "DSt"
$ws.Range("B16").Value = 192

$ws.Range("A17").Value = "DSrt"
$ws.Range("B17").Value = 50

$ws.Range("A18").Value = "DSr"
$ws.Range("B18").Value = 9

$ws.Range("A19").Value = "DSl"
$ws.Range("B19").Value = 193

$ws.Range("A20").Value = "CEtl"
$ws.Range("B20").Value = 418

$ws.Range("A21").Value = "CEr"
$ws.Range("B21").Value = 418

$ws.Range("A22").Value = "CEl"
$ws.Range("B22").Value = 180

$ws.Range("A23").Value = "BWtl"
$ws.Range("B23").Value = 139

$ws.Range("A24").Value = "BWrt"
$ws.Range("B24").Value = 138

$ws.Range("A25").Value = "BWr"
$ws.Range("B25").Value = 177

$ws.Range("A26").Value = "BWl"
$ws.Range("B26").Value = 99

$ws.Range("A27").Value = "ANrt"
$ws.Range("B27").Value = 728

$ws.Range("A28").Value = "ANr"
$ws.Range("B28").Value = 107
